# Apply the commit "[calculateDescriptives] Fix descriptives for pooled variables"
# to the example_varInfo workbook.
#
# Summary of the change:
#  - E3 (IDSCH row) gets the missing German label "Schul-ID".
#  - The three Likert scale-item rows (7-9) get numbered labels
#    "Likert-Skalenindikator 1/2/3" instead of all sharing the same text.
#  - The block that used to be named "pvkat_*" (rows 17-22) is renamed to
#    "pvord_*" (ordinal / Kompetenzstufe plausible values) - this was
#    always meant to represent the ordinal/"Kompetenzstufe" variant.
#  - Row 22 (imputation 5 of that block) had a copy/paste bug where its
#    LabelSH still said "IMPUTATION 4" - fixed to "IMPUTATION 5".
#  - A brand-new block of 6 rows (23-28) is appended for the *nominal*
#    ("Kompetenzkategorie") pooled/plausible-value variables, which is
#    the new "pvkat_*" block (pvkat_pooled, pvkat_1..pvkat_5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: IDSCH - add missing German label in column E -------------------
$ws.Range("E3").Value = "Schul-ID"

# --- Rows 7-9: number the Likert scale indicator labels ---------------------
$ws.Range("E7").Value = "Likert-Skalenindikator 1"
$ws.Range("E8").Value = "Likert-Skalenindikator 2"
$ws.Range("E9").Value = "Likert-Skalenindikator 3"

# --- Rows 17-22: rename the "pvkat_*" variables to "pvord_*" ----------------
# (they represent the ordinal / "Kompetenzstufe" plausible values)
$ws.Range("A17").Value = "pvord_pooled"
$ws.Range("A18").Value = "pvord_1"
$ws.Range("A19").Value = "pvord_2"
$ws.Range("A20").Value = "pvord_3"
$ws.Range("A21").Value = "pvord_4"
$ws.Range("A22").Value = "pvord_5"

# Row 22 LabelSH had a copy/paste bug - it repeated "IMPUTATION 4" instead of "IMPUTATION 5"
$ws.Range("E22").Value = "IMPUTATION 5: Kompetenzstufe des plausible value"

# --- Rows 23-28 (new): the nominal "pvkat_*" plausible-value block ----------
# Row 23: pvkat_pooled
$ws.Range("A23").Value = "pvkat_pooled"
$ws.Range("B23").Value = "sh"
$ws.Range("C23").Value = "PVs"
$ws.Range("D23").Value = "-"
$ws.Range("F23").Value = "-"
$ws.Range("G23").Value = "2.1"
$ws.Range("I23").Value = "nominal plausible value"
$ws.Range("J23").Value = "nein"
$ws.Range("K23").Value = "-"
$ws.Range("L23").Value = "-"
$ws.Range("M23").Value = "nein"
$ws.Range("N23").Value = "-"
$ws.Range("O23").Value = "-"
$ws.Range("P23").Value = "-"
$ws.Range("Q23").Value = "nein"

# Row 24: pvkat_1
$ws.Range("A24").Value = "pvkat_1"
$ws.Range("B24").Value = "ds"
$ws.Range("C24").Value = "PVs"
$ws.Range("D24").Value = "-"
$ws.Range("E24").Value = "IMPUTATION 1: Kompetenzkategorie des plausible value"
$ws.Range("F24").Value = "-"
$ws.Range("G24").Value = "2.1"
$ws.Range("I24").Value = "-"
$ws.Range("J24").Value = "nein"
$ws.Range("K24").Value = "-"
$ws.Range("L24").Value = "-"
$ws.Range("M24").Value = "nein"
$ws.Range("N24").Value = "-"
$ws.Range("O24").Value = "-"
$ws.Range("P24").Value = "-"
$ws.Range("Q24").Value = "nein"

# Row 25: pvkat_2
$ws.Range("A25").Value = "pvkat_2"
$ws.Range("B25").Value = "ds"
$ws.Range("C25").Value = "PVs"
$ws.Range("D25").Value = "-"
$ws.Range("E25").Value = "IMPUTATION 2: Kompetenzkategorie des plausible value"
$ws.Range("F25").Value = "-"
$ws.Range("G25").Value = "2.1"
$ws.Range("I25").Value = "-"
$ws.Range("J25").Value = "nein"
$ws.Range("K25").Value = "-"
$ws.Range("L25").Value = "-"
$ws.Range("M25").Value = "nein"
$ws.Range("N25").Value = "-"
$ws.Range("O25").Value = "-"
$ws.Range("P25").Value = "-"
$ws.Range("Q25").Value = "nein"

# Row 26: pvkat_3
$ws.Range("A26").Value = "pvkat_3"
$ws.Range("B26").Value = "ds"
$ws.Range("C26").Value = "PVs"
$ws.Range("D26").Value = "-"
$ws.Range("E26").Value = "IMPUTATION 3: Kompetenzkategorie des plausible value"
$ws.Range("F26").Value = "-"
$ws.Range("G26").Value = "2.1"
$ws.Range("I26").Value = "-"
$ws.Range("J26").Value = "nein"
$ws.Range("K26").Value = "-"
$ws.Range("L26").Value = "-"
$ws.Range("M26").Value = "nein"
$ws.Range("N26").Value = "-"
$ws.Range("O26").Value = "-"
$ws.Range("P26").Value = "-"
$ws.Range("Q26").Value = "nein"

# Row 27: pvkat_4
$ws.Range("A27").Value = "pvkat_4"
$ws.Range("B27").Value = "ds"
$ws.Range("C27").Value = "PVs"
$ws.Range("D27").Value = "-"
$ws.Range("E27").Value = "IMPUTATION 4: Kompetenzkategorie des plausible value"
$ws.Range("F27").Value = "-"
$ws.Range("G27").Value = "2.1"
$ws.Range("I27").Value = "-"
$ws.Range("J27").Value = "nein"
$ws.Range("K27").Value = "-"
$ws.Range("L27").Value = "-"
$ws.Range("M27").Value = "nein"
$ws.Range("N27").Value = "-"
$ws.Range("O27").Value = "-"
$ws.Range("P27").Value = "-"
$ws.Range("Q27").Value = "nein"

# Row 28: pvkat_5
$ws.Range("A28").Value = "pvkat_5"
$ws.Range("B28").Value = "ds"
$ws.Range("C28").Value = "PVs"
$ws.Range("D28").Value = "-"
$ws.Range("E28").Value = "IMPUTATION 5: Kompetenzkategorie des plausible value"
$ws.Range("F28").Value = "-"
$ws.Range("G28").Value = "2.1"
$ws.Range("I28").Value = "-"
$ws.Range("J28").Value = "nein"
$ws.Range("K28").Value = "-"
$ws.Range("L28").Value = "-"
$ws.Range("M28").Value = "nein"
$ws.Range("N28").Value = "-"
$ws.Range("O28").Value = "-"
$ws.Range("P28").Value = "-"
$ws.Range("Q28").Value = "nein"

Write-Host "Done applying edits."
